$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E ("Kuvaus (FI)" etc. shift right to F:H)
$ws.Columns.Item(5).Insert()

# New header + value for the inserted "Kokonaispisteet" (total points) column
$ws.Range("E4").Value = "Kokonaispisteet"
$ws.Range("E5").Value = "55.44"

# Extend the existing list data validation (previously only D5) to cover the new E5 cell too,
# re-creating it so the merged sqref "D5:E5" matches a single validation rule.
$ws.Range("D5:E5").Validation.Delete()
$ws.Range("D5:E5").Validation.Add(3, 1, 1, """Määrittelemätön,Hyväksyttävissä,Hylätty""")
$ws.Range("D5:E5").Validation.ShowInput = $false
$ws.Range("D5:E5").Validation.ShowError = $false

# Match the saved cursor position recorded in the workbook
$ws.Range("F17").Select() | Out-Null
